$d = $word.ActiveDocument
$d.Content.Find.Execute("Select a user and invite him to a challenge", $true, $false, $false, $false, $false, $true, 1, $false, "Type the user's username and invite him to a challenge", 2)
